$d = $word.ActiveDocument

function Replace-Paragraph([string]$old, [string]$new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        throw "Find.Execute failed to locate target text: $old"
    }
}

$oldP1 = @'
I will be using two sets of data about the different Pokemon you can catch in Pokemon Go. One data set is organized and relatively clean while the other needs to be organized. To make this process easier I will refer to the datasets by the name of the web site I go the off of, so the organized data set will be called Kaggle and the unorganized data will be called RapidAPI. For the Kaggle dataset I will need to read the CSV file into a pandas data frame. The RapidAPI dataset will be much harder, since it is stored in an API I will need to run an API call on it but that is just the start. Since the data is stored in many different APIs I will have to run multiple API calls. I will have to create a for loop to call the different types of APIs which the data is stored in. I will also have to modify the for loop if the information is stored in a different way. 
'@
$newP1 = @'
I used two sets of data about the different types Pokémon you can catch in Pokémon Go. One data set is relatively organized while the other needs to be transformed more. To make this process easier I referred to the datasets by the name of the web site I go the off of, so the organized data set will be called Kaggle and the unorganized data will be called RapidAPI. For the Kaggle dataset I read the CSV file into a pandas data frame. The RapidAPI dataset was much harder, since it was stored in an API and I had to run multiple API calls on it to get the information I needed. I created a for loop to call the different types of APIs which the data is stored in. 
'@

$oldP2 = @'
The plan is to take the Kaggle dataset, rename columns, and review the rows. If the data looks well organized, I will then pull in the Rapid API, select the columns needed, renames the columns, merge tables, remove rows(duplicates), and complete a final review on the data. I will then compare the two datasets and see if there are discrepancies. I will have to do this many times over for the Rapid API dataset, there will be many different panda data frames I will have to search through to get all the information I will need. 
'@
$newP2 = @'
I took the Kaggle dataset in a pandas data and started with making another column to sum up total stats, or total strength of each Pokémon. I then made the legendary column a Boolean to simplify data. For the RapidAPI I started with a definition that pulled the API, grabbed the API key, and put it all into a data frame. I then create a variable that held all the different API calls that I wanted and passed it through a for loop and stored all the data frames into a list of data frames. I then renamed columns and removed duplicate rows, the Pokémon have different forms, and unless the stats differed, I removed the other forms. I then organized the tables reworked the rows and created a combined data frame with all the columns I wanted. I then did the following to both the Kaggle and RapidAPI data sets.  I reordered the columns and split the data frame into three data frames; name, catch, and battle. I then removed duplicates from the catch and battle data frames and created a primary key. Next, I ran a for loop to place a foreign key in the name data frame from the catch and battle data frames. I then checked all the data frames and then created a connection to pgadmin and pushed the data onto the tables.
'@

$oldP3 = @'
When I am finally ready to load the two datasets to SQL I will need to create the PokemonGo database. Then I will have to create two empty tables(RapidAPI and Kaggle) where I will push my data. Next I will create a user and password to allow the python to push data to the dataframe.  Finally I will need to connect to the local database then using pandas load my data frames into the two tables on the PokemonGo database. 
'@
$newP3 = @'
I then created six tables in a pgAdmin SQL database, there is a name, catch, and battle table for the RapidAPI dataset and Kaggle dataset. I chose Pgadmin because I thought it would hold the information best. All the columns needed to be the same and there needed to be values for each row and column. With the stricter rules of pgAdmin this made it very easy. All that was difficult was creating the foreign keys and knowing that the different joins would work. 
'@

Replace-Paragraph $oldP1 $newP1
Replace-Paragraph $oldP2 $newP2
Replace-Paragraph $oldP3 $newP3

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
